# Update cryptocurrency price (column D) and volume-change (column E) figures
# to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "58.199.53"
$ws.Cells.Item(2, 5).Value = "  -0.07%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.597.02"
$ws.Cells.Item(3, 5).Value = "  -0.52%  "
$ws.Cells.Item(4, 5).Value = "  +0.10%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "522.83"
$ws.Cells.Item(5, 5).Value = "  +0.55%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "143.74"
$ws.Cells.Item(6, 5).Value = "  +0.55%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.997"
$ws.Cells.Item(7, 5).Value = "  -0.20%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.569"
$ws.Cells.Item(8, 5).Value = "  +0.43%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "2.617.44"
$ws.Cells.Item(9, 5).Value = "  +0.06%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "6.63"
$ws.Cells.Item(10, 5).Value = "  -0.56%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.101"
$ws.Cells.Item(11, 5).Value = "  -1.57%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.338"
$ws.Cells.Item(12, 5).Value = "  -0.10%  "
$ws.Cells.Item(13, 5).Value = "  -0.39%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "3.058.23"
$ws.Cells.Item(14, 5).Value = "  -0.30%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "58.187.71"
$ws.Cells.Item(15, 5).Value = "  -0.05%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "20.50"
$ws.Cells.Item(16, 5).Value = "  -2.24%  "
$ws.Cells.Item(17, 5).Value = "  -1.41%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "2.585.39"
$ws.Cells.Item(18, 5).Value = "  -0.93%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "339.40"
$ws.Cells.Item(19, 5).Value = "  +1.17%  "
$ws.Cells.Item(20, 5).Value = "  -0.66%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "10.29"
$ws.Cells.Item(21, 5).Value = "  -0.81%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.40"
$ws.Cells.Item(22, 5).Value = "  +1.93%  "
$ws.Cells.Item(23, 5).Value = "  +0.02%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "65.31"
$ws.Cells.Item(24, 5).Value = "  +1.15%  "
$ws.Cells.Item(25, 5).Value = "  +0.47%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.403"
$ws.Cells.Item(26, 5).Value = "  -2.89%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.716.75"
$ws.Cells.Item(27, 5).Value = "  -0.36%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.997"
$ws.Cells.Item(28, 5).Value = "  -0.14%  "
$ws.Cells.Item(29, 5).Value = "  -1.54%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0₃0748"
$ws.Cells.Item(30, 5).Value = "  -5.46%  "
$ws.Cells.Item(31, 5).Value = "  -0.06%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "6.23"
$ws.Cells.Item(32, 5).Value = "  -5.94%  "
$ws.Cells.Item(33, 5).Value = "  -0.02%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "18.80"
$ws.Cells.Item(34, 5).Value = "  +0.21%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "149.71"
$ws.Cells.Item(35, 5).Value = "  -0.24%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "4.03"
$ws.Cells.Item(36, 5).Value = "  -1.70%  "
$ws.Cells.Item(37, 5).Value = "  -3.88%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.875"
$ws.Cells.Item(38, 5).Value = "  -1.63%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.860"
$ws.Cells.Item(39, 5).Value = "  +0.84%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "36.03"
$ws.Cells.Item(40, 5).Value = "  -0.64%  "
$ws.Cells.Item(41, 5).Value = "  +1.68%  "
$ws.Cells.Item(42, 5).Value = "  -2.38%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.996"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "273.20"
$ws.Cells.Item(44, 5).Value = "  +1.32%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.601"
$ws.Cells.Item(45, 5).Value = "  +0.18%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0959"
$ws.Cells.Item(46, 5).Value = "  -0.74%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "10.67"
$ws.Cells.Item(47, 5).Value = "  +0.53%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "18.84"
$ws.Cells.Item(48, 5).Value = "  -1.59%  "
$ws.Cells.Item(49, 5).Value = "  -1.66%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "18.97"
$ws.Cells.Item(50, 5).Value = "  +3.95%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.978.15"
$ws.Cells.Item(51, 5).Value = "  -2.85%  "
